# Update leve-profit calculation columns (H:N) across multiple profession sheets.
# Values come from refreshed Universalis market-price data; only numeric
# cells in columns H-N change (A-G leve metadata is untouched).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3724.2727
$ws.Range("I98").Value = 3421.5
$ws.Range("J98").Value = 6752
$ws.Range("K98").Value = 3421.5
$ws.Range("L98").Value = 6752
$ws.Range("M98").Value = -1923.5
$ws.Range("N98").Value = -9748
$ws.Range("H122").Value = 3724.2727
$ws.Range("I122").Value = 3421.5
$ws.Range("J122").Value = 6752
$ws.Range("K122").Value = 10264.5
$ws.Range("L122").Value = 20256
$ws.Range("M122").Value = -7814.5
$ws.Range("N122").Value = -25156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 13177.429
$ws.Range("I28").Value = 7040.1665
$ws.Range("J28").Value = 50001
$ws.Range("K28").Value = 7040.1665
$ws.Range("L28").Value = 50001
$ws.Range("M28").Value = -6848.1665
$ws.Range("N28").Value = -50385
$ws.Range("H37").Value = 7858.909
$ws.Range("I37").Value = 1278
$ws.Range("K37").Value = 1278
$ws.Range("M37").Value = -1005
$ws.Range("H74").Value = 7888.8696
$ws.Range("I74").Value = 1495.7693
$ws.Range("J74").Value = 16199.9
$ws.Range("K74").Value = 1495.7693
$ws.Range("L74").Value = 16199.9
$ws.Range("M74").Value = -621.7692999999999
$ws.Range("N74").Value = -17947.9
$ws.Range("H77").Value = 7888.8696
$ws.Range("I77").Value = 1495.7693
$ws.Range("J77").Value = 16199.9
$ws.Range("K77").Value = 7478.8465
$ws.Range("L77").Value = 80999.5
$ws.Range("M77").Value = -3110.8465
$ws.Range("N77").Value = -89735.5
$ws.Range("H99").Value = 13177.429
$ws.Range("I99").Value = 7040.1665
$ws.Range("J99").Value = 50001
$ws.Range("K99").Value = 7040.1665
$ws.Range("L99").Value = 50001
$ws.Range("M99").Value = -4045.1665
$ws.Range("N99").Value = -55991
$ws.Range("H122").Value = 1993.5454
$ws.Range("I122").Value = 1704.4615
$ws.Range("J122").Value = 2411.111
$ws.Range("K122").Value = 5113.3845
$ws.Range("L122").Value = 7233.333
$ws.Range("M122").Value = -2663.3845
$ws.Range("N122").Value = -12133.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1450.5
$ws.Range("I86").Value = 1445.1111
$ws.Range("J86").Value = 1466.6666
$ws.Range("K86").Value = 1445.1111
$ws.Range("L86").Value = 1466.6666
$ws.Range("M86").Value = -322.1111000000001
$ws.Range("N86").Value = -3712.6666
$ws.Range("H89").Value = 1450.5
$ws.Range("I89").Value = 1445.1111
$ws.Range("J89").Value = 1466.6666
$ws.Range("K89").Value = 7225.5555
$ws.Range("L89").Value = 7333.333000000001
$ws.Range("M89").Value = -1609.5555
$ws.Range("N89").Value = -18565.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4067928
$ws.Range("I31").Value = 1408.4667
$ws.Range("J31").Value = 6413997
$ws.Range("K31").Value = 1408.4667
$ws.Range("L31").Value = 6413997
$ws.Range("M31").Value = -1113.4667
$ws.Range("N31").Value = -6414587
$ws.Range("H34").Value = 4067928
$ws.Range("I34").Value = 1408.4667
$ws.Range("J34").Value = 6413997
$ws.Range("K34").Value = 1408.4667
$ws.Range("L34").Value = 6413997
$ws.Range("M34").Value = -1206.4667
$ws.Range("N34").Value = -6414401
$ws.Range("H50").Value = 10535.2
$ws.Range("J50").Value = 10535.2
$ws.Range("L50").Value = 10535.2
$ws.Range("N50").Value = -11785.2
$ws.Range("H51").Value = 8768.700000000001
$ws.Range("I51").Value = 5095
$ws.Range("J51").Value = 9687.125
$ws.Range("K51").Value = 5095
$ws.Range("L51").Value = 9687.125
$ws.Range("M51").Value = -4359
$ws.Range("N51").Value = -11159.125
$ws.Range("H59").Value = 14508.2
$ws.Range("I59").Value = 10500
$ws.Range("J59").Value = 15510.25
$ws.Range("K59").Value = 10500
$ws.Range("L59").Value = 15510.25
$ws.Range("M59").Value = -9355
$ws.Range("N59").Value = -17800.25
$ws.Range("H60").Value = 8276
$ws.Range("I60").Value = 3000
$ws.Range("J60").Value = 8755.637000000001
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 8755.637000000001
$ws.Range("M60").Value = -2489
$ws.Range("N60").Value = -9777.637000000001
$ws.Range("H61").Value = 8768.700000000001
$ws.Range("I61").Value = 5095
$ws.Range("J61").Value = 9687.125
$ws.Range("K61").Value = 5095
$ws.Range("L61").Value = 9687.125
$ws.Range("M61").Value = -4747
$ws.Range("N61").Value = -10383.125
$ws.Range("H62").Value = 2875
$ws.Range("I62").Value = 2667.9167
$ws.Range("J62").Value = 3151.111
$ws.Range("K62").Value = 2667.9167
$ws.Range("L62").Value = 3151.111
$ws.Range("M62").Value = -2043.9167
$ws.Range("N62").Value = -4399.111
$ws.Range("H65").Value = 2875
$ws.Range("I65").Value = 2667.9167
$ws.Range("J65").Value = 3151.111
$ws.Range("K65").Value = 13339.5835
$ws.Range("L65").Value = 15755.555
$ws.Range("M65").Value = -10219.5835
$ws.Range("N65").Value = -21995.555
$ws.Range("H68").Value = 18306.666
$ws.Range("J68").Value = 19114.4
$ws.Range("L68").Value = 19114.4
$ws.Range("N68").Value = -20612.4
$ws.Range("H71").Value = 18306.666
$ws.Range("J71").Value = 19114.4
$ws.Range("L71").Value = 57343.2
$ws.Range("N71").Value = -64831.2
$ws.Range("H97").Value = 13900
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H122").Value = 3575838
$ws.Range("I122").Value = 4634745.5
$ws.Range("J122").Value = 2025
$ws.Range("K122").Value = 13904236.5
$ws.Range("L122").Value = 6075
$ws.Range("M122").Value = -13901786.5
$ws.Range("N122").Value = -10975
$ws.Range("H134").Value = 2020.7059
$ws.Range("I134").Value = 1888.825
$ws.Range("J134").Value = 2500.2727
$ws.Range("K134").Value = 5666.475
$ws.Range("L134").Value = 7500.8181
$ws.Range("M134").Value = -3131.475
$ws.Range("N134").Value = -12570.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 3197.5
$ws.Range("J35").Value = 3930
$ws.Range("L35").Value = 11790
$ws.Range("N35").Value = -12366
$ws.Range("H70").Value = 37096596
$ws.Range("I70").Value = 37096596
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 111289788
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -111289473
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 37096596
$ws.Range("I73").Value = 37096596
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 111289788
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -111288696
$ws.Range("N73").ClearContents()
$ws.Range("H75").Value = 428.75
$ws.Range("I75").Value = 100
$ws.Range("J75").Value = 538.3333
$ws.Range("K75").Value = 300
$ws.Range("L75").Value = 1614.9999
$ws.Range("M75").Value = 698
$ws.Range("N75").Value = -3610.9999
$ws.Range("H78").Value = 428.75
$ws.Range("I78").Value = 100
$ws.Range("J78").Value = 538.3333
$ws.Range("K78").Value = 900
$ws.Range("L78").Value = 4844.9997
$ws.Range("M78").Value = 4092
$ws.Range("N78").Value = -14828.9997
$ws.Range("H103").Value = 2429359.8
$ws.Range("I103").Value = 5667137.5
$ws.Range("J103").Value = 1026.625
$ws.Range("K103").Value = 17001412.5
$ws.Range("L103").Value = 3079.875
$ws.Range("M103").Value = -17000533.5
$ws.Range("N103").Value = -4837.875
$ws.Range("H113").Value = 781.4286
$ws.Range("I113").Value = 615.2308
$ws.Range("J113").Value = 879.63635
$ws.Range("K113").Value = 1845.6924
$ws.Range("L113").Value = 2638.90905
$ws.Range("M113").Value = 324.3075999999999
$ws.Range("N113").Value = -6978.90905
$ws.Range("H131").Value = 2624.0168
$ws.Range("J131").Value = 1688.5178
$ws.Range("L131").Value = 5065.553400000001
$ws.Range("N131").Value = -15145.5534

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1965.625
$ws.Range("I97").Value = 1795.4546
$ws.Range("J97").Value = 2340
$ws.Range("K97").Value = 1795.4546
$ws.Range("L97").Value = 2340
$ws.Range("M97").Value = -1299.4546
$ws.Range("N97").Value = -3332
$ws.Range("H107").Value = 420.85715
$ws.Range("I107").Value = 296.06668
$ws.Range("J107").Value = 732.8333
$ws.Range("K107").Value = 296.06668
$ws.Range("L107").Value = 732.8333
$ws.Range("M107").Value = 1623.93332
$ws.Range("N107").Value = -4572.8333
$ws.Range("H122").Value = 3168
$ws.Range("I122").Value = 3119.2
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 9357.599999999999
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -6907.599999999999
$ws.Range("N122").Value = -16600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1727
$ws.Range("I40").Value = 1469.3334
$ws.Range("K40").Value = 1469.3334
$ws.Range("M40").Value = -1333.3334
$ws.Range("H99").Value = 20000
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 2052.6428
$ws.Range("I122").Value = 1891.069
$ws.Range("J122").Value = 2413.077
$ws.Range("K122").Value = 5673.207
$ws.Range("L122").Value = 7239.231000000001
$ws.Range("M122").Value = -3223.207
$ws.Range("N122").Value = -12139.231
$ws.Range("H136").Value = 2108.2666
$ws.Range("I136").Value = 1134.9166
$ws.Range("J136").Value = 6001.6665
$ws.Range("K136").Value = 3404.7498
$ws.Range("L136").Value = 18004.9995
$ws.Range("M136").Value = -854.7498000000001
$ws.Range("N136").Value = -23104.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4712.75
$ws.Range("J62").Value = 4475
$ws.Range("L62").Value = 4475
$ws.Range("N62").Value = -5723
$ws.Range("H65").Value = 4712.75
$ws.Range("J65").Value = 4475
$ws.Range("L65").Value = 22375
$ws.Range("N65").Value = -28615
$ws.Range("H96").Value = 1552.1
$ws.Range("I96").Value = 1639.5454
$ws.Range("J96").Value = 1445.2222
$ws.Range("K96").Value = 1639.5454
$ws.Range("L96").Value = 1445.2222
$ws.Range("M96").Value = -266.5454
$ws.Range("N96").Value = -4191.2222
$ws.Range("H122").Value = 2169.238
$ws.Range("I122").Value = 1581.7142
$ws.Range("J122").Value = 3344.2856
$ws.Range("K122").Value = 4745.142599999999
$ws.Range("L122").Value = 10032.8568
$ws.Range("M122").Value = -2295.142599999999
$ws.Range("N122").Value = -14932.8568
$ws.Range("H132").Value = 10604473
$ws.Range("I132").Value = 12260984
$ws.Range("J132").Value = 2799.6
$ws.Range("K132").Value = 36782952
$ws.Range("L132").Value = 8398.799999999999
$ws.Range("M132").Value = -36780422
$ws.Range("N132").Value = -13458.8
